# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) used emoji characters (📕 📘 📙 📗) as status
# markers. These are replaced with plain-text / simple-symbol equivalents:
#   📕 -> -3
#   📘 -> ⚠️
#   📙 -> +3
#   📗 -> ✅

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old emoji value -> new value
$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📙" = "+3"
    "📗" = "✅"
}

# Replacement values that look like numbers must be forced to stay text,
# otherwise Excel would store them as numeric cells instead of strings.
$forceText = @{
    "-3" = $true
    "+3" = $true
}

$lastRow = $ws.Cells(1, 1).SpecialCells(11).Row  # xlCellTypeLastCell
$col = 1  # column A = "statut"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Text

    if ($map.ContainsKey($current)) {
        $newValue = $map[$current]

        if ($forceText.ContainsKey($newValue)) {
            # Force text interpretation so "-3"/"+3" aren't turned into numbers
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            # Restore the default "Normal" style so no extra formatting lingers
            $cell.Style = "Normal"
        }
        else {
            $cell.Value = $newValue
        }
    }
}
